# Estadisticos Matutinos 15 Oct
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Estadisticos 1P" (sheet1.xml): update Blancos(D)/Aprobados(F)/
# Por_Apro(G) for rows 2-5 and add Promedio(H) values.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")

$ws1.Range("D2").Value = 11
$ws1.Range("F2").Value = 20
$ws1.Range("G2").Value = 64.52
$ws1.Range("H2").Value = 8.5

$ws1.Range("D3").Value = 10
$ws1.Range("F3").Value = 11
$ws1.Range("G3").Value = 52.38
$ws1.Range("H3").Value = 9.300000000000001

$ws1.Range("D4").Value = 9
$ws1.Range("F4").Value = 26
$ws1.Range("G4").Value = 74.29000000000001
$ws1.Range("H4").Value = 8.5

$ws1.Range("D5").Value = 13
$ws1.Range("F5").Value = 20
$ws1.Range("G5").Value = 60.61
$ws1.Range("H5").Value = 8.199999999999999

# ---------------------------------------------------------------------------
# Sheet "Estadisticos 2P" (sheet2.xml): only Reprobados(E) changes.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")

$ws2.Range("E2").Value = 20
$ws2.Range("E3").Value = 11
$ws2.Range("E4").Value = 26
$ws2.Range("E5").Value = 20

# ---------------------------------------------------------------------------
# Sheet "Estadisticos Final" (sheet3.xml): same pattern as "Estadisticos 1P".
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")

$ws3.Range("D2").Value = 11
$ws3.Range("F2").Value = 20
$ws3.Range("G2").Value = 64.52
$ws3.Range("H2").Value = 8.5

$ws3.Range("D3").Value = 10
$ws3.Range("F3").Value = 11
$ws3.Range("G3").Value = 52.38
$ws3.Range("H3").Value = 9.300000000000001

$ws3.Range("D4").Value = 9
$ws3.Range("F4").Value = 26
$ws3.Range("G4").Value = 74.29000000000001
$ws3.Range("H4").Value = 8.5

$ws3.Range("D5").Value = 13
$ws3.Range("F5").Value = 20
$ws3.Range("G5").Value = 60.61
$ws3.Range("H5").Value = 8.199999999999999

# ---------------------------------------------------------------------------
# Sheet "Rescatables" (sheet4.xml): append six new student rows (2-7).
# Columns: A=NC, B=Paterno, C=Materno, D=Nombres, E=Nombre_Largo, F=Grupo,
# G=Reprobadas
# Values are written column-by-column (matching how the source data was
# pasted in) so new shared-string entries land in the same order as the
# original edit: all Paterno, then all Materno, then all Nombres, ...
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")

$ws4.Range("A2").Value = 20330051920136
$ws4.Range("A3").Value = 20330051920151
$ws4.Range("A4").Value = 19330051920407
$ws4.Range("A5").Value = 20330051920191
$ws4.Range("A6").Value = 20330051920202
$ws4.Range("A7").Value = 20330051920270

$ws4.Range("B2").Value = "MARCIAL"
$ws4.Range("B3").Value = "TORRES"
$ws4.Range("B4").Value = "ROSAS"
$ws4.Range("B5").Value = "ARIAS"
$ws4.Range("B6").Value = "DE JESUS"
$ws4.Range("B7").Value = "HERNANDEZ"

$ws4.Range("C2").Value = "MORALES"
$ws4.Range("C3").Value = "PEREZ"
$ws4.Range("C4").Value = "SANCHEZ"
$ws4.Range("C5").Value = "BARRAGAN"
$ws4.Range("C6").Value = "CASTILLO"
$ws4.Range("C7").Value = "DIAZ"

$ws4.Range("D2").Value = "IVAN DE JESUS"
$ws4.Range("D3").Value = "CONSTANZA XIMENA"
$ws4.Range("D4").Value = "MONSERRAT"
$ws4.Range("D5").Value = "ANALI"
$ws4.Range("D6").Value = "ITZEL"
$ws4.Range("D7").Value = "ELIAS ALONSO"

$ws4.Range("E2").Value = "ÉTICA"
$ws4.Range("E3").Value = "ÉTICA"
$ws4.Range("E4").Value = "CIENCIA, TECNOLOGÍA, SOCIEDAD Y VALORES"
$ws4.Range("E5").Value = "ÉTICA"
$ws4.Range("E6").Value = "ÉTICA"
$ws4.Range("E7").Value = "ÉTICA"

$ws4.Range("F2").Value = "3ARHV"
$ws4.Range("F3").Value = "3ARHV"
$ws4.Range("F4").Value = "5ASV"
$ws4.Range("F5").Value = "3ALCV"
$ws4.Range("F6").Value = "3ALCV"
$ws4.Range("F7").Value = "3APV"

$ws4.Range("G2").Value = 6
$ws4.Range("G3").Value = 6
$ws4.Range("G4").Value = 6
$ws4.Range("G5").Value = 6
$ws4.Range("G6").Value = 6
$ws4.Range("G7").Value = 6
